# Auto-generated script applying numeric cell updates per xml diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 260.7
$ws.Range("K6").Value = 322.5
$ws.Range("I6").Value = 107.5
$ws.Range("M6").Value = -210.5
$ws.Range("I18").Value = 4500
$ws.Range("K18").Value = 4500
$ws.Range("M18").Value = -4216
$ws.Range("H18").Value = 6000
$ws.Range("I33").Value = 1244.55
$ws.Range("H33").Value = 992
$ws.Range("K33").Value = 1244.55
$ws.Range("M33").Value = -1015.55
$ws.Range("K53").Value = 542.25
$ws.Range("N53").Value = -4274.6667
$ws.Range("H53").Value = 1595.8572
$ws.Range("I53").Value = 542.25
$ws.Range("J53").Value = 3000.6667
$ws.Range("M53").Value = 94.75
$ws.Range("L53").Value = 3000.6667
$ws.Range("M64").Value = -7080.3335
$ws.Range("I64").Value = 7328.3335
$ws.Range("K64").Value = 7328.3335
$ws.Range("H64").Value = 8996.25
$ws.Range("K67").Value = 7328.3335
$ws.Range("M67").Value = -6470.3335
$ws.Range("H67").Value = 8996.25
$ws.Range("I67").Value = 7328.3335
$ws.Range("L80").Value = 9990.900000000001
$ws.Range("N80").Value = -11986.9
$ws.Range("I80").Value = 1054.3125
$ws.Range("K80").Value = 3162.9375
$ws.Range("J80").Value = 3330.3
$ws.Range("M80").Value = -2164.9375
$ws.Range("M82").Value = -12194
$ws.Range("H82").Value = 4200
$ws.Range("K82").Value = 12600
$ws.Range("I82").Value = 4200
$ws.Range("M83").Value = -4496.8125
$ws.Range("I83").Value = 1054.3125
$ws.Range("J83").Value = 3330.3
$ws.Range("N83").Value = -39956.7
$ws.Range("K83").Value = 9488.8125
$ws.Range("L83").Value = 29972.7
$ws.Range("H85").Value = 4200
$ws.Range("M85").Value = -11196
$ws.Range("K85").Value = 12600
$ws.Range("I85").Value = 4200
$ws.Range("M86").Value = 121.5
$ws.Range("K86").Value = 1001.5
$ws.Range("I86").Value = 1001.5
$ws.Range("H86").Value = 4002.5
$ws.Range("I89").Value = 1001.5
$ws.Range("H89").Value = 4002.5
$ws.Range("M89").Value = 608.5
$ws.Range("N92").Value = -2826
$ws.Range("K92").Value = 282.81818
$ws.Range("M92").Value = 965.18182
$ws.Range("L92").Value = 330
$ws.Range("J92").Value = 330
$ws.Range("I92").Value = 282.81818
$ws.Range("H92").Value = 286.75
$ws.Range("K98").Value = 4159.385
$ws.Range("J98").Value = 6003
$ws.Range("M98").Value = -2661.385
$ws.Range("L98").Value = 6003
$ws.Range("H98").Value = 4505.0625
$ws.Range("N98").Value = -8999
$ws.Range("I98").Value = 4159.385
$ws.Range("I100").Value = 1177.5834
$ws.Range("K100").Value = 1177.5834
$ws.Range("H100").Value = 6304.4165
$ws.Range("M100").Value = -636.5834
$ws.Range("N113").Value = -16346.5
$ws.Range("J113").Value = 9838.5
$ws.Range("H113").Value = 7878.75
$ws.Range("L113").Value = 9838.5
$ws.Range("K122").Value = 12478.155
$ws.Range("N122").Value = -22909
$ws.Range("J122").Value = 6003
$ws.Range("L122").Value = 18009
$ws.Range("I122").Value = 4159.385
$ws.Range("M122").Value = -10028.155
$ws.Range("H122").Value = 4505.0625
$ws.Range("L125").Value = 51426
$ws.Range("N125").Value = -56346
$ws.Range("K125").Value = 17532
$ws.Range("I125").Value = 1948
$ws.Range("M125").Value = -15072
$ws.Range("H125").Value = 4877.1113
$ws.Range("J125").Value = 5714
$ws.Range("J129").Value = 1711.4231
$ws.Range("N129").Value = -15134.2693
$ws.Range("L129").Value = 5134.2693
$ws.Range("H129").Value = 1312.3846
$ws.Range("H132").Value = 6596.1
$ws.Range("M132").Value = -17187.3038
$ws.Range("I132").Value = 6572.4346
$ws.Range("K132").Value = 19717.3038
$ws.Range("I135").Value = 8375.583000000001
$ws.Range("K135").Value = 75380.247
$ws.Range("M135").Value = -72845.247
$ws.Range("H135").Value = 10521.682
$ws.Range("I137").Value = 3041.9092
$ws.Range("H137").Value = 2964.6765
$ws.Range("M137").Value = -6575.7276
$ws.Range("K137").Value = 9125.7276
$ws.Range("M138").Value = 1536.1177
$ws.Range("H138").Value = 1621895.4
$ws.Range("K138").Value = 3603.8823
$ws.Range("J138").Value = 2723967.5
$ws.Range("N138").Value = -8182182.5
$ws.Range("I138").Value = 1201.2941
$ws.Range("L138").Value = 8171902.5
$ws.Range("J141").Value = 3165.5
$ws.Range("N141").Value = -19856.5
$ws.Range("L141").Value = 9496.5
$ws.Range("H141").Value = 2336.4517
$ws.Range("I141").Value = 2137.48
$ws.Range("M141").Value = -1232.440000000001
$ws.Range("K141").Value = 6412.440000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M5").Value = -315.16666
$ws.Range("H5").Value = 427.16666
$ws.Range("I5").Value = 427.16666
$ws.Range("K5").Value = 427.16666
$ws.Range("L41").Value = 9950
$ws.Range("J41").Value = 9950
$ws.Range("N41").Value = -10778
$ws.Range("H41").Value = 5500.5
$ws.Range("I61").Value = 5723.2583
$ws.Range("L61").Value = 10610.333
$ws.Range("K61").Value = 5723.2583
$ws.Range("J61").Value = 10610.333
$ws.Range("N61").Value = -11034.333
$ws.Range("H61").Value = 6515.757
$ws.Range("M61").Value = -5511.2583
$ws.Range("M74").Value = -4340
$ws.Range("K74").Value = 5214
$ws.Range("I74").Value = 5214
$ws.Range("H74").Value = 6139.2
$ws.Range("I77").Value = 5214
$ws.Range("M77").Value = -21702
$ws.Range("H77").Value = 6139.2
$ws.Range("K77").Value = 26070
$ws.Range("J88").Value = 2878.7778
$ws.Range("K88").Value = 2125
$ws.Range("N88").Value = -3690.7778
$ws.Range("I88").Value = 2125
$ws.Range("H88").Value = 2690.3333
$ws.Range("M88").Value = -1719
$ws.Range("L88").Value = 2878.7778
$ws.Range("N91").Value = -5686.7778
$ws.Range("K91").Value = 2125
$ws.Range("M91").Value = -721
$ws.Range("I91").Value = 2125
$ws.Range("L91").Value = 2878.7778
$ws.Range("H91").Value = 2690.3333
$ws.Range("J91").Value = 2878.7778
$ws.Range("M102").Value = 554.75
$ws.Range("L102").Value = 40001510
$ws.Range("H102").Value = 15385854
$ws.Range("K102").Value = 1067.25
$ws.Range("N102").Value = -40004754
$ws.Range("I102").Value = 1067.25
$ws.Range("J102").Value = 40001510
$ws.Range("K122").Value = 9281.143199999999
$ws.Range("I122").Value = 3093.7144
$ws.Range("M122").Value = -6831.143199999999
$ws.Range("H122").Value = 3035
$ws.Range("J129").Value = 80000
$ws.Range("N129").Value = -90000
$ws.Range("L129").Value = 80000
$ws.Range("H129").Value = 80000
$ws.Range("H132").Value = 2873.3015
$ws.Range("J132").Value = 5156.7144
$ws.Range("M132").Value = -4132.694
$ws.Range("N132").Value = -20530.1432
$ws.Range("I132").Value = 2220.898
$ws.Range("K132").Value = 6662.694
$ws.Range("L132").Value = 15470.1432
$ws.Range("J136").Value = 10610.333
$ws.Range("I136").Value = 5723.2583
$ws.Range("L136").Value = 31830.999
$ws.Range("N136").Value = -36930.999
$ws.Range("M136").Value = -14619.7749
$ws.Range("H136").Value = 6515.757
$ws.Range("K136").Value = 17169.7749

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M4").Value = -312.16666
$ws.Range("K4").Value = 427.16666
$ws.Range("H4").Value = 427.16666
$ws.Range("I4").Value = 427.16666
$ws.Range("N20").Value = -2243.5
$ws.Range("H20").Value = 1568.8572
$ws.Range("J20").Value = 1749.5
$ws.Range("L20").Value = 1749.5
$ws.Range("K22").Value = 462
$ws.Range("M22").Value = -289
$ws.Range("H22").Value = 482.6875
$ws.Range("I22").Value = 462
$ws.Range("H27").Value = 80000
$ws.Range("L27").Value = 80000
$ws.Range("J27").Value = 80000
$ws.Range("N27").Value = -80384
$ws.Range("M86").Value = -66124.25
$ws.Range("K86").Value = 67247.25
$ws.Range("L86").Value = 9395.799999999999
$ws.Range("N86").Value = -11641.8
$ws.Range("I86").Value = 67247.25
$ws.Range("H86").Value = 35107.555
$ws.Range("J86").Value = 9395.799999999999
$ws.Range("I89").Value = 67247.25
$ws.Range("H89").Value = 35107.555
$ws.Range("K89").Value = 336236.25
$ws.Range("L89").Value = 46979
$ws.Range("J89").Value = 9395.799999999999
$ws.Range("M89").Value = -330620.25
$ws.Range("N89").Value = -58211
$ws.Range("M96").Value = -11728
$ws.Range("I96").Value = 14474
$ws.Range("H96").Value = 14474
$ws.Range("K96").Value = 14474
$ws.Range("H99").Value = 2900.6667
$ws.Range("J99").Value = 4700.5
$ws.Range("L99").Value = 4700.5
$ws.Range("M99").Value = -502.75
$ws.Range("N99").Value = -7696.5
$ws.Range("K99").Value = 2000.75
$ws.Range("I99").Value = 2000.75
$ws.Range("K105").Value = 2037.1282
$ws.Range("H105").Value = 2956.541
$ws.Range("I105").Value = 2037.1282
$ws.Range("M105").Value = -290.1282000000001
$ws.Range("L107").Value = 380
$ws.Range("N107").Value = -4220
$ws.Range("H107").Value = 2569.8823
$ws.Range("J107").Value = 380
$ws.Range("I134").Value = 2970.16
$ws.Range("H134").Value = 2970.16
$ws.Range("K134").Value = 8910.48
$ws.Range("M134").Value = -6375.48

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I31").Value = 4237.0625
$ws.Range("K31").Value = 4237.0625
$ws.Range("H31").Value = 6212.7334
$ws.Range("M31").Value = -3942.0625
$ws.Range("K34").Value = 4237.0625
$ws.Range("I34").Value = 4237.0625
$ws.Range("M34").Value = -4035.0625
$ws.Range("H34").Value = 6212.7334
$ws.Range("H58").Value = 1446.5405
$ws.Range("L58").Value = 1333.1818
$ws.Range("J58").Value = 1333.1818
$ws.Range("M58").Value = -1291.5
$ws.Range("I58").Value = 1494.5
$ws.Range("N58").Value = -1739.1818
$ws.Range("K58").Value = 1494.5
$ws.Range("M86").Value = -5219
$ws.Range("K86").Value = 6342
$ws.Range("I86").Value = 6342
$ws.Range("H86").Value = 11033.846
$ws.Range("I89").Value = 6342
$ws.Range("H89").Value = 11033.846
$ws.Range("K89").Value = 31710
$ws.Range("M89").Value = -26094
$ws.Range("H99").Value = 3975.8823
$ws.Range("J99").Value = 3156.875
$ws.Range("L99").Value = 3156.875
$ws.Range("M99").Value = -3205.8887
$ws.Range("N99").Value = -6152.875
$ws.Range("K99").Value = 4703.8887
$ws.Range("I99").Value = 4703.8887
$ws.Range("K105").Value = 3212
$ws.Range("H105").Value = 3268.7144
$ws.Range("I105").Value = 3212
$ws.Range("M105").Value = -1465
$ws.Range("L126").Value = 9470.625
$ws.Range("M126").Value = -11641.6661
$ws.Range("J126").Value = 3156.875
$ws.Range("H126").Value = 3975.8823
$ws.Range("N126").Value = -14410.625
$ws.Range("I126").Value = 4703.8887
$ws.Range("K126").Value = 14111.6661
$ws.Range("N134").Value = -12468
$ws.Range("I134").Value = 1996.5
$ws.Range("H134").Value = 2278.2
$ws.Range("J134").Value = 2466
$ws.Range("K134").Value = 5989.5
$ws.Range("M134").Value = -3454.5
$ws.Range("L134").Value = 7398
$ws.Range("J136").Value = 1333.1818
$ws.Range("I136").Value = 1494.5
$ws.Range("L136").Value = 3999.5454
$ws.Range("N136").Value = -9099.545399999999
$ws.Range("M136").Value = -1933.5
$ws.Range("H136").Value = 1446.5405
$ws.Range("K136").Value = 4483.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 270
$ws.Range("N2").Value = -1848.00004
$ws.Range("L2").Value = 1622.00004
$ws.Range("J2").Value = 270.33334
$ws.Range("L5").Value = 4174.9998
$ws.Range("J5").Value = 1391.6666
$ws.Range("H5").Value = 1391.6666
$ws.Range("N5").Value = -4398.9998
$ws.Range("N33").Value = -776
$ws.Range("L33").Value = 210
$ws.Range("J33").Value = 35
$ws.Range("H33").Value = 39.083332
$ws.Range("L60").Value = 9004.5
$ws.Range("N60").Value = -9506.5
$ws.Range("H60").Value = 2141
$ws.Range("J60").Value = 3001.5
$ws.Range("I62").Value = 3011
$ws.Range("L62").Value = 30813.999
$ws.Range("M62").Value = -8347
$ws.Range("J62").Value = 10271.333
$ws.Range("K62").Value = 9033
$ws.Range("H62").Value = 7367.2
$ws.Range("N62").Value = -32185.999
$ws.Range("J63").Value = 7994
$ws.Range("N63").Value = -25480
$ws.Range("L63").Value = 23982
$ws.Range("H63").Value = 4448.25
$ws.Range("M64").Value = -6396
$ws.Range("I64").Value = 2222
$ws.Range("K64").Value = 6666
$ws.Range("H64").Value = 4611
$ws.Range("N65").Value = -99305.997
$ws.Range("I65").Value = 3011
$ws.Range("K65").Value = 27099
$ws.Range("L65").Value = 92441.997
$ws.Range("H65").Value = 7367.2
$ws.Range("J65").Value = 10271.333
$ws.Range("M65").Value = -23667
$ws.Range("N66").Value = -79434
$ws.Range("L66").Value = 71946
$ws.Range("J66").Value = 7994
$ws.Range("H66").Value = 4448.25
$ws.Range("K67").Value = 6666
$ws.Range("M67").Value = -5730
$ws.Range("H67").Value = 4611
$ws.Range("I67").Value = 2222
$ws.Range("L107").Value = 3439.8
$ws.Range("N107").Value = -7279.799999999999
$ws.Range("H107").Value = 1059.5
$ws.Range("J107").Value = 1146.6
$ws.Range("H108").Value = 1128.5
$ws.Range("K108").Value = 3385.5
$ws.Range("M108").Value = -505.5
$ws.Range("I108").Value = 1128.5
$ws.Range("I109").Value = 3259.889
$ws.Range("J109").Value = 7798.6
$ws.Range("K109").Value = 9779.667000000001
$ws.Range("M109").Value = -8739.667000000001
$ws.Range("N109").Value = -25475.8
$ws.Range("H109").Value = 5648.684
$ws.Range("L109").Value = 23395.8
$ws.Range("N113").Value = -8044.25
$ws.Range("J113").Value = 1234.75
$ws.Range("H113").Value = 1056.5625
$ws.Range("L113").Value = 3704.25
$ws.Range("J124").Value = 10226.5
$ws.Range("H124").Value = 6446.2
$ws.Range("N124").Value = -40499.5
$ws.Range("L124").Value = 30679.5
$ws.Range("J129").Value = 6511.923
$ws.Range("N129").Value = -29535.769
$ws.Range("L129").Value = 19535.769
$ws.Range("H129").Value = 3950.5217
$ws.Range("N130").Value = -37256
$ws.Range("H130").Value = 3375.8572
$ws.Range("I130").Value = 1097.4
$ws.Range("K130").Value = 3292.2
$ws.Range("J130").Value = 9072
$ws.Range("M130").Value = 1727.8
$ws.Range("L130").Value = 27216
$ws.Range("M131").Value = -27484.5
$ws.Range("K131").Value = 32524.5
$ws.Range("L131").Value = 7330.0905
$ws.Range("I131").Value = 10841.5
$ws.Range("N131").Value = -17410.0905
$ws.Range("H131").Value = 4082.0244
$ws.Range("J131").Value = 2443.3635
$ws.Range("H132").Value = 3400
$ws.Range("J132").Value = 3400
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -35660
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 30600
$ws.Range("N135").Value = -17594.9994
$ws.Range("L135").Value = 12524.9994
$ws.Range("H135").Value = 1391.6666
$ws.Range("J135").Value = 1391.6666
$ws.Range("L140").Value = 12281.7276
$ws.Range("N140").Value = -22641.7276
$ws.Range("H140").Value = 3402.2742
$ws.Range("J140").Value = 4093.9092

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 3621.6
$ws.Range("K6").Value = 702.6667
$ws.Range("I6").Value = 702.6667
$ws.Range("J6").Value = 8000
$ws.Range("M6").Value = -589.6667
$ws.Range("N6").Value = -8226
$ws.Range("L6").Value = 8000
$ws.Range("M7").Value = -15389
$ws.Range("H7").Value = 263306
$ws.Range("I7").Value = 15501
$ws.Range("K7").Value = 15501
$ws.Range("I8").Value = 15501
$ws.Range("H8").Value = 263306
$ws.Range("K8").Value = 15501
$ws.Range("M8").Value = -15362
$ws.Range("H11").Value = 3999999.8
$ws.Range("L11").Value = 4999999.5
$ws.Range("J11").Value = 4999999.5
$ws.Range("N11").Value = -5000277.5
$ws.Range("J16").Value = 8000
$ws.Range("I16").Value = 702.6667
$ws.Range("H16").Value = 3621.6
$ws.Range("M16").Value = -452.6667
$ws.Range("L16").Value = 8000
$ws.Range("K16").Value = 702.6667
$ws.Range("N16").Value = -8500
$ws.Range("N93").Value = -91006.75
$ws.Range("J93").Value = 87262.75
$ws.Range("L93").Value = 87262.75
$ws.Range("H93").Value = 87262.75
$ws.Range("M102").Value = 109.8888999999999
$ws.Range("H102").Value = 49057.906
$ws.Range("K102").Value = 1512.1111
$ws.Range("I102").Value = 1512.1111
$ws.Range("J106").Value = 48177.4
$ws.Range("L106").Value = 48177.4
$ws.Range("N106").Value = -50701.4
$ws.Range("H106").Value = 48177.4
$ws.Range("K122").Value = 300321690
$ws.Range("I122").Value = 100107230
$ws.Range("M122").Value = -300319240
$ws.Range("H122").Value = 83424240
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -509777.99
$ws.Range("J126").Value = 0
$ws.Range("H126").Value = 170749.33
$ws.Range("N126").ClearContents()
$ws.Range("I126").Value = 170749.33
$ws.Range("K126").Value = 512247.99
$ws.Range("H132").Value = 3598.0308
$ws.Range("M132").Value = -10114.4552
$ws.Range("I132").Value = 4214.8184
$ws.Range("K132").Value = 12644.4552

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M7").Value = -5638
$ws.Range("H7").Value = 7027.6665
$ws.Range("N7").Value = -11723.5
$ws.Range("I7").Value = 5750
$ws.Range("L7").Value = 11499.5
$ws.Range("J7").Value = 11499.5
$ws.Range("K7").Value = 5750
$ws.Range("L22").Value = 9375
$ws.Range("J22").Value = 9375
$ws.Range("H22").Value = 8877.888999999999
$ws.Range("N22").Value = -9965
$ws.Range("H27").Value = 8877.888999999999
$ws.Range("L27").Value = 9375
$ws.Range("J27").Value = 9375
$ws.Range("N27").Value = -9589
$ws.Range("I61").Value = 7022.636
$ws.Range("L61").Value = 5333.3335
$ws.Range("K61").Value = 7022.636
$ws.Range("J61").Value = 5333.3335
$ws.Range("N61").Value = -5737.3335
$ws.Range("H61").Value = 6819.92
$ws.Range("M61").Value = -6820.636
$ws.Range("M93").Value = 822.2777699999999
$ws.Range("I93").Value = 425.72223
$ws.Range("H93").Value = 657.7619
$ws.Range("K93").Value = 425.72223
$ws.Range("I100").Value = 3521.7778
$ws.Range("K100").Value = 3521.7778
$ws.Range("H100").Value = 3846.3333
$ws.Range("M100").Value = -2980.7778
$ws.Range("I113").Value = 7022.636
$ws.Range("N113").Value = -9673.333500000001
$ws.Range("K113").Value = 7022.636
$ws.Range("J113").Value = 5333.3335
$ws.Range("H113").Value = 6819.92
$ws.Range("L113").Value = 5333.3335
$ws.Range("M113").Value = -4852.636
$ws.Range("L126").Value = 34498.5
$ws.Range("M126").Value = -14780
$ws.Range("J126").Value = 11499.5
$ws.Range("H126").Value = 7027.6665
$ws.Range("N126").Value = -39438.5
$ws.Range("I126").Value = 5750
$ws.Range("K126").Value = 17250
$ws.Range("H132").Value = 7977.325
$ws.Range("J132").Value = 8764.25
$ws.Range("M132").Value = -19828.1255
$ws.Range("N132").Value = -31352.75
$ws.Range("I132").Value = 7452.7085
$ws.Range("K132").Value = 22358.1255
$ws.Range("L132").Value = 26292.75
$ws.Range("J136").Value = 4165.8335
$ws.Range("I136").Value = 5139.241
$ws.Range("L136").Value = 12497.5005
$ws.Range("N136").Value = -17597.5005
$ws.Range("M136").Value = -12867.723
$ws.Range("H136").Value = 4972.3716
$ws.Range("K136").Value = 15417.723

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K9").Value = 5006
$ws.Range("I9").Value = 5006
$ws.Range("H9").Value = 7003
$ws.Range("M9").Value = -4866
$ws.Range("N9").Value = -9280
$ws.Range("J9").Value = 9000
$ws.Range("L9").Value = 9000
$ws.Range("H14").Value = 4999
$ws.Range("I14").Value = 6000
$ws.Range("L14").Value = 995
$ws.Range("J14").Value = 995
$ws.Range("N14").Value = -1331
$ws.Range("K14").Value = 6000
$ws.Range("M14").Value = -5832
$ws.Range("J96").Value = 3724.5
$ws.Range("H96").Value = 2644.3635
$ws.Range("N96").Value = -6470.5
$ws.Range("L96").Value = 3724.5
$ws.Range("M107").Value = -354
$ws.Range("L107").Value = 2684.6667
$ws.Range("N107").Value = -6524.6667
$ws.Range("I107").Value = 758
$ws.Range("H107").Value = 822.8421
$ws.Range("J107").Value = 894.8889
$ws.Range("K107").Value = 2274
$ws.Range("K122").Value = 10079.625
$ws.Range("N122").Value = -10900
$ws.Range("J122").Value = 2000
$ws.Range("L122").Value = 6000
$ws.Range("I122").Value = 3359.875
$ws.Range("M122").Value = -7629.625
$ws.Range("H122").Value = 3305.48
$ws.Range("L126").Value = 14755.8
$ws.Range("M126").Value = -6787.000100000001
$ws.Range("J126").Value = 4918.6
$ws.Range("H126").Value = 3347.5144
$ws.Range("N126").Value = -19695.8
$ws.Range("I126").Value = 3085.6667
$ws.Range("K126").Value = 9257.000100000001
$ws.Range("H132").Value = 1393.3684
$ws.Range("M132").Value = -1451
$ws.Range("I132").Value = 1327
$ws.Range("K132").Value = 3981
$ws.Range("J136").Value = 6201.1763
$ws.Range("I136").Value = 5157.971
$ws.Range("L136").Value = 18603.5289
$ws.Range("N136").Value = -23703.5289
$ws.Range("M136").Value = -12923.913
$ws.Range("H136").Value = 5499.019
$ws.Range("K136").Value = 15473.913

